$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D3 formula to use absolute reference for A1
$ws.Range("D3").Formula = "=A$1*1.1"

# D4:D34 become a shared formula block referencing A$1*1.1
$ws.Range("D4:D34").Formula = "=A$1*1.1"

# D35 becomes a standalone formula referencing D3
$ws.Range("D35").Formula = "=D3/3"

# D36:D38 become a shared formula block referencing D4
$ws.Range("D36:D38").Formula = "=D4/3"

# Update the sheet view: scroll position and selection
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D35:D38").Select()
